$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data rows (2-8) get cyclically re-ordered within two groups:
#   Group A: rows 2,3,4,5  -> new row r takes the old values of row r+1 (row 5 takes old row 2)
#   Group B: rows 6,7,8    -> new row r takes the old values of row r+1 (row 8 takes old row 6)
# Columns A,B,C,E,F,G,H,I,J,K are identical for every row, so only D,L,M,N,O,P,Q,R,S,T need updating.

function Get-RowData($ws, $r) {
    return @{
        D = $ws.Cells.Item($r, 4).Value2
        L = $ws.Cells.Item($r, 12).Value2
        M = $ws.Cells.Item($r, 13).Value2
        N = $ws.Cells.Item($r, 14).Value2
        O = $ws.Cells.Item($r, 15).Value2
        P = $ws.Cells.Item($r, 16).Value2
        Q = $ws.Cells.Item($r, 17).Value2
        R = $ws.Cells.Item($r, 18).Value2
        S = $ws.Cells.Item($r, 19).Value2
        T = $ws.Cells.Item($r, 20).Value2
    }
}

function Set-RowData($ws, $r, $data) {
    $ws.Cells.Item($r, 4).Value = $data.D
    $ws.Cells.Item($r, 12).Value = $data.L
    $ws.Cells.Item($r, 13).Value = $data.M
    $ws.Cells.Item($r, 14).Value = $data.N
    $ws.Cells.Item($r, 15).Value = $data.O
    $ws.Cells.Item($r, 16).Value = $data.P
    $ws.Cells.Item($r, 17).Value = $data.Q
    $ws.Cells.Item($r, 18).Value = $data.R
    $ws.Cells.Item($r, 19).Value = $data.S
    $ws.Cells.Item($r, 20).Value = $data.T
}

# Capture original values before overwriting anything
$row2 = Get-RowData $ws 2
$row3 = Get-RowData $ws 3
$row4 = Get-RowData $ws 4
$row5 = Get-RowData $ws 5
$row6 = Get-RowData $ws 6
$row7 = Get-RowData $ws 7
$row8 = Get-RowData $ws 8

# Apply the cyclic shift
Set-RowData $ws 2 $row3
Set-RowData $ws 3 $row4
Set-RowData $ws 4 $row5
Set-RowData $ws 5 $row2

Set-RowData $ws 6 $row7
Set-RowData $ws 7 $row8
Set-RowData $ws 8 $row6
